# Scheduled runner update: refresh Kraken market-price snapshots (H/I/J/K/L)
# and recomputed leve profit figures (M/N) across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 2000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""

$ws.Range("H32").Value = 10999.857
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 10999.857
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10999.857
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -11651.857

$ws.Range("H62").Value = 2750
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -2126

$ws.Range("H65").Value = 2750
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -10630

$ws.Range("H80").Value = 7887.5
$ws.Range("I80").Value = 516.6667
$ws.Range("J80").Value = 30000
$ws.Range("K80").Value = 1550.0001
$ws.Range("L80").Value = 90000
$ws.Range("M80").Value = -552.0001
$ws.Range("N80").Value = -91996

$ws.Range("H83").Value = 7887.5
$ws.Range("I83").Value = 516.6667
$ws.Range("J83").Value = 30000
$ws.Range("K83").Value = 4650.0003
$ws.Range("L83").Value = 270000
$ws.Range("M83").Value = 341.9997000000003
$ws.Range("N83").Value = -279984

$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("K86").Value = 1000
$ws.Range("M86").Value = 123

$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("K89").Value = 5000
$ws.Range("M89").Value = 616

$ws.Range("H100").Value = 1722.3334
$ws.Range("I100").Value = 826.625
$ws.Range("J100").Value = 8888
$ws.Range("K100").Value = 826.625
$ws.Range("L100").Value = 8888
$ws.Range("M100").Value = -285.625
$ws.Range("N100").Value = -9970

$ws.Range("H138").Value = 3897.8
$ws.Range("J138").Value = 4996.3335
$ws.Range("L138").Value = 14989.0005
$ws.Range("N138").Value = -25269.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 80000
$ws.Range("I82").Value = 80000
$ws.Range("K82").Value = 80000
$ws.Range("M82").Value = -79639

$ws.Range("H85").Value = 80000
$ws.Range("I85").Value = 80000
$ws.Range("K85").Value = 80000
$ws.Range("M85").Value = -78752

$ws.Range("H97").Value = 565
$ws.Range("I97").Value = 565
$ws.Range("K97").Value = 565
$ws.Range("M97").Value = -69

$ws.Range("H102").Value = 2699.4
$ws.Range("I102").Value = 2699.4
$ws.Range("K102").Value = 2699.4
$ws.Range("M102").Value = -1077.4

$ws.Range("H132").Value = 4478
$ws.Range("I132").Value = 4478
$ws.Range("K132").Value = 13434
$ws.Range("M132").Value = -10904

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 34700
$ws.Range("I26").Value = 34700
$ws.Range("K26").Value = 34700
$ws.Range("M26").Value = -34408

$ws.Range("H99").Value = 2699.8
$ws.Range("I99").Value = 1999.5
$ws.Range("K99").Value = 1999.5
$ws.Range("M99").Value = -501.5

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = ""
$ws.Range("N140").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 878.7143
$ws.Range("I58").Value = 878.7143
$ws.Range("K58").Value = 878.7143
$ws.Range("M58").Value = -675.7143

$ws.Range("H106").Value = 26501.834
$ws.Range("J106").Value = 26501.834
$ws.Range("L106").Value = 26501.834
$ws.Range("N106").Value = -29025.834

$ws.Range("H136").Value = 878.7143
$ws.Range("I136").Value = 878.7143
$ws.Range("K136").Value = 2636.1429
$ws.Range("M136").Value = -86.14289999999983

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 577.8333
$ws.Range("I5").Value = 493.4
$ws.Range("K5").Value = 1480.2
$ws.Range("M5").Value = -1368.2

$ws.Range("H62").Value = 4999
$ws.Range("J62").Value = 4999
$ws.Range("L62").Value = 14997
$ws.Range("N62").Value = -16369

$ws.Range("H65").Value = 4999
$ws.Range("J65").Value = 4999
$ws.Range("L65").Value = 44991
$ws.Range("N65").Value = -51855

$ws.Range("H113").Value = 580.3570999999999
$ws.Range("J113").Value = 629.7273
$ws.Range("L113").Value = 1889.1819
$ws.Range("N113").Value = -6229.1819

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""

$ws.Range("H135").Value = 577.8333
$ws.Range("I135").Value = 493.4
$ws.Range("K135").Value = 4440.599999999999
$ws.Range("M135").Value = -1905.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 417416.5
$ws.Range("J3").Value = 833
$ws.Range("L3").Value = 833
$ws.Range("N3").Value = -1065

$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002

$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008

$ws.Range("H132").Value = 5736.364
$ws.Range("I132").Value = 5122.222
$ws.Range("K132").Value = 15366.666
$ws.Range("M132").Value = -12836.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = 0

$ws.Range("H61").Value = 7824.75
$ws.Range("I61").Value = 7824.75
$ws.Range("K61").Value = 7824.75
$ws.Range("M61").Value = -7622.75

$ws.Range("H93").Value = 40000
$ws.Range("I93").Value = 40000
$ws.Range("K93").Value = 40000
$ws.Range("M93").Value = -38752

$ws.Range("H100").Value = 4167.6665
$ws.Range("I100").Value = 3003
$ws.Range("K100").Value = 3003
$ws.Range("M100").Value = -2462

$ws.Range("H113").Value = 7824.75
$ws.Range("I113").Value = 7824.75
$ws.Range("K113").Value = 7824.75
$ws.Range("M113").Value = -5654.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 1010001.8
$ws.Range("I24").Value = 1668269.6
$ws.Range("J24").Value = 22600
$ws.Range("K24").Value = 1668269.6
$ws.Range("L24").Value = 22600
$ws.Range("M24").Value = -1668039.6
$ws.Range("N24").Value = -23060

$ws.Range("H54").Value = 1500
$ws.Range("I54").Value = 1500
$ws.Range("K54").Value = 1500
$ws.Range("M54").Value = -980

$ws.Range("H58").Value = 40000
$ws.Range("I58").Value = 40000
$ws.Range("K58").Value = 40000
$ws.Range("M58").Value = -39692

$ws.Range("H96").Value = 2500
$ws.Range("I96").Value = 2500
$ws.Range("K96").Value = 2500
$ws.Range("M96").Value = -1127

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = ""
$ws.Range("N105").Value = 0
